$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 362.46155
$ws.Range("J12").Value = 766.3333
$ws.Range("L12").Value = 766.3333
$ws.Range("N12").Value = -1106.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1738
$ws.Range("J19").Value = 2182
$ws.Range("L19").Value = 2182
$ws.Range("N19").Value = -2532

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 5568
$ws.Range("I31").Value = 4566.6665
$ws.Range("K31").Value = 13699.9995
$ws.Range("M31").Value = -13469.9995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2929.6155
$ws.Range("I70").Value = 3999
$ws.Range("J70").Value = 2735.182
$ws.Range("K70").Value = 11997
$ws.Range("L70").Value = 8205.545999999998
$ws.Range("M70").Value = -11727
$ws.Range("N70").Value = -8745.545999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2929.6155
$ws.Range("I73").Value = 3999
$ws.Range("J73").Value = 2735.182
$ws.Range("K73").Value = 11997
$ws.Range("L73").Value = 8205.545999999998
$ws.Range("M73").Value = -11061
$ws.Range("N73").Value = -10077.546

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 16016.952
$ws.Range("I76").Value = 8763.200000000001
$ws.Range("J76").Value = 16997.19
$ws.Range("K76").Value = 8763.200000000001
$ws.Range("L76").Value = 16997.19
$ws.Range("M76").Value = -8448.200000000001
$ws.Range("N76").Value = -17627.19

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 16016.952
$ws.Range("I79").Value = 8763.200000000001
$ws.Range("J79").Value = 16997.19
$ws.Range("K79").Value = 8763.200000000001
$ws.Range("L79").Value = 16997.19
$ws.Range("M79").Value = -7671.200000000001
$ws.Range("N79").Value = -19181.19

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 827.5625
$ws.Range("I80").Value = 664.3333
$ws.Range("K80").Value = 1992.9999
$ws.Range("M80").Value = -994.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 827.5625
$ws.Range("I83").Value = 664.3333
$ws.Range("K83").Value = 5978.9997
$ws.Range("M83").Value = -986.9997000000003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 4799.8
$ws.Range("J88").Value = 4874.75
$ws.Range("L88").Value = 4874.75
$ws.Range("N88").Value = -5686.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 4799.8
$ws.Range("J91").Value = 4874.75
$ws.Range("L91").Value = 4874.75
$ws.Range("N91").Value = -7682.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1705.25
$ws.Range("J129").Value = 1808.5
$ws.Range("L129").Value = 5425.5
$ws.Range("N129").Value = -15425.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2728.48
$ws.Range("I138").Value = 1785.875
$ws.Range("J138").Value = 3172.0588
$ws.Range("K138").Value = 5357.625
$ws.Range("L138").Value = 9516.1764
$ws.Range("M138").Value = -217.625
$ws.Range("N138").Value = -19796.1764

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2960.8118
$ws.Range("I32").Value = 2960.8118
$ws.Range("K32").Value = 2960.8118
$ws.Range("M32").Value = -2673.8118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 12816.827
$ws.Range("I61").Value = 7940.409
$ws.Range("K61").Value = 7940.409
$ws.Range("M61").Value = -7728.409

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2627
$ws.Range("I122").Value = 1885
$ws.Range("J122").Value = 2998
$ws.Range("K122").Value = 5655
$ws.Range("L122").Value = 8994
$ws.Range("M122").Value = -3205
$ws.Range("N122").Value = -13894

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 124999
$ws.Range("J134").Value = 124999
$ws.Range("L134").Value = 124999
$ws.Range("N134").Value = -135139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 12816.827
$ws.Range("I136").Value = 7940.409
$ws.Range("K136").Value = 23821.227
$ws.Range("M136").Value = -21271.227

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 199900
$ws.Range("J43").Value = 199900
$ws.Range("L43").Value = 199900
$ws.Range("N43").Value = -200262

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2525.2334
$ws.Range("I134").Value = 2496.724
$ws.Range("K134").Value = 7490.172
$ws.Range("M134").Value = -4955.172

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 57.409092
$ws.Range("I7").Value = 52
$ws.Range("J7").Value = 81.75
$ws.Range("K7").Value = 52
$ws.Range("L7").Value = 81.75
$ws.Range("M7").Value = 61
$ws.Range("N7").Value = -307.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 661.625
$ws.Range("I22").Value = 235.33333
$ws.Range("K22").Value = 235.33333
$ws.Range("M22").Value = 114.66667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7825.923
$ws.Range("I31").Value = 6249.25
$ws.Range("K31").Value = 6249.25
$ws.Range("M31").Value = -5954.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7825.923
$ws.Range("I34").Value = 6249.25
$ws.Range("K34").Value = 6249.25
$ws.Range("M34").Value = -6047.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 49000
$ws.Range("I98").Value = 44000
$ws.Range("K98").Value = 44000
$ws.Range("M98").Value = -41754

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5675.5386
$ws.Range("I132").Value = 6116.5454
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 18349.6362
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -15819.6362
$ws.Range("N132").Value = -14810

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9347
$ws.Range("I134").Value = 7651.2
$ws.Range("J134").Value = 14999.667
$ws.Range("K134").Value = 22953.6
$ws.Range("L134").Value = 44999.001
$ws.Range("M134").Value = -20418.6
$ws.Range("N134").Value = -50069.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 175.16667
$ws.Range("J23").Value = 250
$ws.Range("L23").Value = 750
$ws.Range("N23").Value = -1220

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 870.4286
$ws.Range("J38").Value = 1021.5
$ws.Range("L38").Value = 3064.5
$ws.Range("N38").Value = -3758.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 379
$ws.Range("I107").Value = 249.25
$ws.Range("J107").Value = 453.14285
$ws.Range("K107").Value = 747.75
$ws.Range("L107").Value = 1359.42855
$ws.Range("M107").Value = 1172.25
$ws.Range("N107").Value = -5199.428550000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 639.5
$ws.Range("I118").Value = 639.5
$ws.Range("K118").Value = 1918.5
$ws.Range("M118").Value = -675.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 685.8421
$ws.Range("J122").Value = 816.46155
$ws.Range("L122").Value = 7348.15395
$ws.Range("N122").Value = -12248.15395

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 42546.832
$ws.Range("J136").Value = 42546.832
$ws.Range("L136").Value = 127640.496
$ws.Range("N136").Value = -132740.496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1514.4546
$ws.Range("J22").Value = 1507.5
$ws.Range("L22").Value = 1507.5
$ws.Range("N22").Value = -2097.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1514.4546
$ws.Range("J27").Value = 1507.5
$ws.Range("L27").Value = 1507.5
$ws.Range("N27").Value = -1721.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 4788
$ws.Range("I35").Value = 4788
$ws.Range("K35").Value = 4788
$ws.Range("M35").Value = -4452

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6463.3955
$ws.Range("I93").Value = 1680.8667
$ws.Range("K93").Value = 1680.8667
$ws.Range("M93").Value = -432.8667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4049.5715
$ws.Range("I122").Value = 3799.4546
$ws.Range("J122").Value = 4966.6665
$ws.Range("K122").Value = 11398.3638
$ws.Range("L122").Value = 14899.9995
$ws.Range("M122").Value = -8948.363799999999
$ws.Range("N122").Value = -19799.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3454.889
$ws.Range("I136").Value = 3267.6829
$ws.Range("K136").Value = 9803.048699999999
$ws.Range("M136").Value = -7253.048699999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6037.25
$ws.Range("I132").Value = 6037.25
$ws.Range("K132").Value = 18111.75
$ws.Range("M132").Value = -15581.75
